$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the staff edit spinner value: Hours Worked for row 2 (G2) goes from 16 to 17
$ws.Range("G2").Value = 17.0

# Update the active selection left behind on the sheet after the edit
[void]$ws.Range("M10").Select()
